# "Signed Off time sheets"
# Fill in the Supervisor Name field and the second (supervisor) sign-off
# row on the weekly timesheet: supervisor name, initials and sign-off date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supervisor Name: (merged G6:I6)
$ws.Range("G6").Value = "Ankita Gangotra"

# Second sign-off row (merged A27:C27 / D27:E27): initials + date signed
$ws.Range("A27").Value = "A.G"
$ws.Range("D27").Value = 41800
$ws.Range("D27").NumberFormat = "mm-dd-yy"

# Leave the selection on the cell that was last edited, like Excel would
$ws.Range("D27:E27").Select() | Out-Null
